$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B68: convert inline string "1" to a real numeric value 1
$ws.Range("B68").Value = 1

# Add new row 69 with annotation data
$ws.Range("A69").Value = "Ying Tang"
$ws.Range("B69").Value = "'2"
$ws.Range("C69").Value = " weak"
$ws.Range("D69").Value = "CRT"
$ws.Range("E69").Value = "RES"
$ws.Range("F69").Value = "6325282a-75f6-4567-8bb3-3102657c705c"
$ws.Range("G69").Value = "fm5jfAwPbOfP6_annotated.xlsx"
$ws.Range("H69").Value = "I found the empirical evaluation to be weak."
